$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '57.132.10'
$ws.Range("E2").Value = '  +4.77%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.328.99'
$ws.Range("E3").Value = '  +1.84%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '518.95'
$ws.Range("E5").Value = '  +3.05%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.61'
$ws.Range("E6").Value = '  +3.39%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.32%  '
$ws.Range("E8").Value = '  +1.98%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.355.17'
$ws.Range("E9").Value = '  +2.39%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.104'
$ws.Range("E10").Value = '  +7.99%  '
$ws.Range("E11").Value = '  +0.88%  '
$ws.Range("E12").Value = '  +6.56%  '
$ws.Range("E13").Value = '  +1.93%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.92'
$ws.Range("E14").Value = '  +3.58%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.777.26'
$ws.Range("E15").Value = '  +3.12%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '56.985.25'
$ws.Range("E16").Value = '  +4.47%  '
$ws.Range("E17").Value = '  +3.67%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.356.22'
$ws.Range("E18").Value = '  +2.70%  '
$ws.Range("E19").Value = '  +2.13%  '
$ws.Range("E20").Value = '  +3.50%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '323.20'
$ws.Range("E21").Value = '  +5.48%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.71'
$ws.Range("E22").Value = '  +5.95%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '61.37'
$ws.Range("E24").Value = '  +0.94%  '
$ws.Range("B25").Value = 'Binance-PegBSC-USD'
$ws.Range("C25").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").Value = '  +0.42%  '
$ws.Range("B26").Value = 'Kaspa'
$ws.Range("C26").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.160'
$ws.Range("E26").Value = '  +6.77%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.80'
$ws.Range("E27").Value = '  +5.58%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '172.01'
$ws.Range("E28").Value = '  -0.28%  '
$ws.Range("B29").Value = 'PEPE'
$ws.Range("C29").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0741'
$ws.Range("E29").Value = '  +4.51%  '
$ws.Range("B30").Value = 'Fetch.AI'
$ws.Range("C30").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.21'
$ws.Range("E30").Value = '  +9.56%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.30'
$ws.Range("E31").Value = '  +3.86%  '
$ws.Range("E32").Value = '  +3.62%  '
$ws.Range("E33").Value = '  +2.47%  '
$ws.Range("E34").Value = '  +0.07%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.958'
$ws.Range("E35").Value = '  +2.42%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.996'
$ws.Range("E36").Value = '  +0.04%  '
$ws.Range("E37").Value = '  +4.66%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.03'
$ws.Range("E38").Value = '  +7.15%  '
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.52'
$ws.Range("E39").Value = '  +7.28%  '
$ws.Range("B40").Value = 'OKB'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '37.57'
$ws.Range("E40").Value = '  +3.78%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.383'
$ws.Range("E41").Value = '  +1.65%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '140.25'
$ws.Range("E42").Value = '  +11.97%  '
$ws.Range("E43").Value = '  +5.49%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '277.22'
$ws.Range("E44").Value = '  +11.90%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.17'
$ws.Range("E45").Value = '  +1.82%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0512'
$ws.Range("E46").Value = '  +3.76%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0932'
$ws.Range("E48").Value = '  +2.54%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.382'
$ws.Range("E49").Value = '  +1.72%  '
$ws.Range("E50").Value = '  +4.44%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '16.97'
$ws.Range("E51").Value = '  +2.56%  '
